# Insert a new data row before current row 147 (below the header row),
# which pushes the existing rows 147-262 down to 148-263, then populate
# the newly-inserted row 147 with the new reading.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 147 - existing row 147 (and everything
# below it) shifts down by one row.
$ws.Rows("147:147").Insert()

# Populate the new row 147 with the new record's values. The columns that
# are constant across the whole dataset (A, B, C, E, F, G, H, Q, R) are
# copied from the neighbouring rows; the rest are the new reading's data.
$ws.Range("A147").Value2 = 3
$ws.Range("B147").Value2 = "Femacal de La Calera"
$ws.Range("C147").Value2 = "Coquimbo"
$ws.Range("D147").Value2 = 44494
$ws.Range("E147").Value2 = 5
$ws.Range("F147").Value2 = 100112028
$ws.Range("G147").Value2 = "Sandia"
$ws.Range("H147").Value2 = "Sin especificar"
$ws.Range("I147").Value2 = "Primera"
$ws.Range("J147").Value2 = 310
$ws.Range("K147").Value2 = 750
$ws.Range("L147").Value2 = 800
$ws.Range("M147").Value2 = 774
$ws.Range("N147").Value2 = "`$/kilo (volumen en unidades)"
$ws.Range("O147").Value2 = "Perú"
$ws.Range("P147").Value2 = 774
$ws.Range("Q147").Value2 = 1
$ws.Range("R147").Value2 = "Hortaliza"

# Make sure the date cell keeps the date/time number format used by the
# rest of the "Fecha" column.
$ws.Range("D147").NumberFormat = $ws.Range("D148").NumberFormat
